# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 6488
    6  = 1946
    7  = 1481
    10 = 351
    12 = 5618
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
